$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells "I0" (col I) and "IF" (col J), matching the style used
# by the other header cells (e.g. H1 "IP") - bold, bordered, centered.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New data values for columns I (I0) and J (IF), rows 2-30. I and J hold
# the same value for every row.
$values = @{
    2  = 8
    3  = 5
    4  = 8
    5  = 7
    6  = 9
    7  = 9
    8  = 7
    9  = 9
    10 = 7
    11 = 9
    12 = 8
    13 = 8
    14 = 6
    15 = 7
    16 = 6
    17 = 7
    18 = 5
    19 = 8
    20 = 6
    21 = 6
    22 = 7
    23 = 9
    24 = 8
    25 = 6
    26 = 6
    27 = 5
    28 = 8
    29 = 9
    30 = 7
}

foreach ($row in $values.Keys) {
    $v = $values[$row]
    $ws.Cells.Item($row, 9).Value = $v
    $ws.Cells.Item($row, 10).Value = $v
}
